# Auto-generated script applying cached value updates to Sheets (scheduled runner refresh)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1236.138
$ws.Range("I28").Value = 820.2273
$ws.Range("J28").Value = 2543.2856
$ws.Range("K28").Value = 820.2273
$ws.Range("L28").Value = 2543.2856
$ws.Range("M28").Value = -335.2273
$ws.Range("N28").Value = -3513.2856
$ws.Range("H33").Value = 449.76923
$ws.Range("I33").Value = 380.7
$ws.Range("K33").Value = 380.7
$ws.Range("M33").Value = -151.7
$ws.Range("H34").Value = 7199.75
$ws.Range("I34").Value = 7199.75
$ws.Range("K34").Value = 7199.75
$ws.Range("M34").Value = -6996.75
$ws.Range("H36").Value = 7199.75
$ws.Range("I36").Value = 7199.75
$ws.Range("K36").Value = 7199.75
$ws.Range("M36").Value = -6484.75
$ws.Range("H43").Value = 7128.625
$ws.Range("I43").Value = 694
$ws.Range("J43").Value = 9273.5
$ws.Range("K43").Value = 694
$ws.Range("L43").Value = 9273.5
$ws.Range("M43").Value = -625
$ws.Range("N43").Value = -9411.5
$ws.Range("H53").Value = 2519.7693
$ws.Range("I53").Value = 526.3333
$ws.Range("K53").Value = 526.3333
$ws.Range("M53").Value = 110.6667
$ws.Range("H92").Value = 3794.7
$ws.Range("I92").Value = 2534.9092
$ws.Range("J92").Value = 5334.4443
$ws.Range("K92").Value = 2534.9092
$ws.Range("L92").Value = 5334.4443
$ws.Range("M92").Value = -1286.9092
$ws.Range("N92").Value = -7830.4443
$ws.Range("H98").Value = 291565.94
$ws.Range("I98").Value = 1391.4667
$ws.Range("K98").Value = 1391.4667
$ws.Range("M98").Value = 106.5333000000001
$ws.Range("H101").Value = 710.1177
$ws.Range("I101").Value = 648.3333
$ws.Range("J101").Value = 779.625
$ws.Range("K101").Value = 1944.9999
$ws.Range("L101").Value = 2338.875
$ws.Range("M101").Value = -322.9999
$ws.Range("N101").Value = -5582.875
$ws.Range("H113").Value = 4798.5293
$ws.Range("I113").Value = 3535.7778
$ws.Range("J113").Value = 6219.125
$ws.Range("K113").Value = 3535.7778
$ws.Range("L113").Value = 6219.125
$ws.Range("M113").Value = -281.7777999999998
$ws.Range("N113").Value = -12727.125
$ws.Range("H122").Value = 291565.94
$ws.Range("I122").Value = 1391.4667
$ws.Range("K122").Value = 4174.4001
$ws.Range("M122").Value = -1724.4001
$ws.Range("H132").Value = 1355
$ws.Range("I132").Value = 1118.8182
$ws.Range("J132").Value = 3953
$ws.Range("K132").Value = 3356.4546
$ws.Range("L132").Value = 11859
$ws.Range("M132").Value = -826.4546
$ws.Range("N132").Value = -16919
$ws.Range("H138").Value = 2637.4075
$ws.Range("I138").Value = 1184.4839
$ws.Range("K138").Value = 3553.4517
$ws.Range("M138").Value = 1586.5483
$ws.Range("H141").Value = 2611.4285
$ws.Range("I141").Value = 2611.4285
$ws.Range("K141").Value = 7834.2855
$ws.Range("M141").Value = -2654.2855

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 529.15
$ws.Range("I2").Value = 383.3125
$ws.Range("J2").Value = 1112.5
$ws.Range("K2").Value = 383.3125
$ws.Range("L2").Value = 1112.5
$ws.Range("M2").Value = -270.3125
$ws.Range("N2").Value = -1338.5
$ws.Range("H61").Value = 4025.8647
$ws.Range("I61").Value = 2929.5312
$ws.Range("K61").Value = 2929.5312
$ws.Range("M61").Value = -2717.5312
$ws.Range("H102").Value = 2218.75
$ws.Range("I102").Value = 2028.5714
$ws.Range("J102").Value = 3550
$ws.Range("K102").Value = 2028.5714
$ws.Range("L102").Value = 3550
$ws.Range("M102").Value = -406.5714
$ws.Range("N102").Value = -6794
$ws.Range("H116").Value = 529.15
$ws.Range("I116").Value = 383.3125
$ws.Range("J116").Value = 1112.5
$ws.Range("K116").Value = 383.3125
$ws.Range("L116").Value = 1112.5
$ws.Range("M116").Value = 1910.6875
$ws.Range("N116").Value = -5700.5
$ws.Range("H122").Value = 3943.8572
$ws.Range("I122").Value = 3777.2
$ws.Range("J122").Value = 4095.3635
$ws.Range("K122").Value = 11331.6
$ws.Range("L122").Value = 12286.0905
$ws.Range("M122").Value = -8881.599999999999
$ws.Range("N122").Value = -17186.0905
$ws.Range("H136").Value = 4025.8647
$ws.Range("I136").Value = 2929.5312
$ws.Range("K136").Value = 8788.5936
$ws.Range("M136").Value = -6238.5936

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 529.15
$ws.Range("I3").Value = 383.3125
$ws.Range("J3").Value = 1112.5
$ws.Range("K3").Value = 383.3125
$ws.Range("L3").Value = 1112.5
$ws.Range("M3").Value = -269.3125
$ws.Range("N3").Value = -1340.5
$ws.Range("H94").Value = 2309.5
$ws.Range("I94").Value = 2304.5
$ws.Range("K94").Value = 2304.5
$ws.Range("M94").Value = -1853.5
$ws.Range("H99").Value = 2759.353
$ws.Range("I99").Value = 2564.0833
$ws.Range("J99").Value = 3228
$ws.Range("K99").Value = 2564.0833
$ws.Range("L99").Value = 3228
$ws.Range("M99").Value = -1066.0833
$ws.Range("N99").Value = -6224
$ws.Range("H134").Value = 2656.95
$ws.Range("I134").Value = 2656.95
$ws.Range("K134").Value = 7970.849999999999
$ws.Range("M134").Value = -5435.849999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 110894
$ws.Range("I31").Value = 9506
$ws.Range("J31").Value = 136241
$ws.Range("K31").Value = 9506
$ws.Range("L31").Value = 136241
$ws.Range("M31").Value = -9211
$ws.Range("N31").Value = -136831
$ws.Range("H34").Value = 110894
$ws.Range("I34").Value = 9506
$ws.Range("J34").Value = 136241
$ws.Range("K34").Value = 9506
$ws.Range("L34").Value = 136241
$ws.Range("M34").Value = -9304
$ws.Range("N34").Value = -136645
$ws.Range("H134").Value = 2983.5881
$ws.Range("I134").Value = 2054.7693
$ws.Range("K134").Value = 6164.3079
$ws.Range("M134").Value = -3629.3079

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H50").Value = 27784508
$ws.Range("I50").Value = 55557380
$ws.Range("K50").Value = 166672140
$ws.Range("M50").Value = -166671659
$ws.Range("H53").Value = 27784508
$ws.Range("I53").Value = 55557380
$ws.Range("K53").Value = 166672140
$ws.Range("M53").Value = -166671659
$ws.Range("H86").Value = 6196.5454
$ws.Range("J86").Value = 6427.5713
$ws.Range("L86").Value = 19282.7139
$ws.Range("N86").Value = -21654.7139
$ws.Range("H89").Value = 6196.5454
$ws.Range("J89").Value = 6427.5713
$ws.Range("L89").Value = 57848.14169999999
$ws.Range("N89").Value = -69704.14169999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2684.5208
$ws.Range("I102").Value = 1910.6
$ws.Range("J102").Value = 4768.154
$ws.Range("K102").Value = 1910.6
$ws.Range("L102").Value = 4768.154
$ws.Range("M102").Value = -288.5999999999999
$ws.Range("N102").Value = -8012.154
$ws.Range("H122").Value = 7686.5806
$ws.Range("I122").Value = 9156.777
$ws.Range("K122").Value = 27470.331
$ws.Range("M122").Value = -25020.331
$ws.Range("H126").Value = 4501.1
$ws.Range("I126").Value = 3160.2727
$ws.Range("J126").Value = 6139.8887
$ws.Range("K126").Value = 9480.8181
$ws.Range("L126").Value = 18419.6661
$ws.Range("M126").Value = -7010.8181
$ws.Range("N126").Value = -23359.6661
$ws.Range("H127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("N127").ClearContents()
$ws.Range("H132").Value = 2328.842
$ws.Range("I132").Value = 1218.4166
$ws.Range("J132").Value = 4232.4287
$ws.Range("K132").Value = 3655.2498
$ws.Range("L132").Value = 12697.2861
$ws.Range("M132").Value = -1125.2498
$ws.Range("N132").Value = -17757.2861

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2774.4285
$ws.Range("J22").Value = 4349.5835
$ws.Range("L22").Value = 4349.5835
$ws.Range("N22").Value = -4939.5835
$ws.Range("H27").Value = 2774.4285
$ws.Range("J27").Value = 4349.5835
$ws.Range("L27").Value = 4349.5835
$ws.Range("N27").Value = -4563.5835
$ws.Range("H46").Value = 3437.2368
$ws.Range("J46").Value = 4458.2173
$ws.Range("L46").Value = 4458.2173
$ws.Range("N46").Value = -4834.2173
$ws.Range("H55").Value = 1713.8928
$ws.Range("I55").Value = 477.4375
$ws.Range("K55").Value = 477.4375
$ws.Range("M55").Value = -304.4375
$ws.Range("H100").Value = 11790.947
$ws.Range("I100").Value = 4332.6665
$ws.Range("J100").Value = 13189.375
$ws.Range("K100").Value = 4332.6665
$ws.Range("L100").Value = 13189.375
$ws.Range("M100").Value = -3791.6665
$ws.Range("N100").Value = -14271.375
$ws.Range("H130").Value = 84856.664
$ws.Range("J130").Value = 84856.664
$ws.Range("L130").Value = 84856.664
$ws.Range("N130").Value = -94896.664
$ws.Range("H132").Value = 5847.4443
$ws.Range("I132").Value = 4786.4287
$ws.Range("J132").Value = 6522.636
$ws.Range("K132").Value = 14359.2861
$ws.Range("L132").Value = 19567.908
$ws.Range("M132").Value = -11829.2861
$ws.Range("N132").Value = -24627.908

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H29").Value = 1533750
$ws.Range("I29").Value = 120000
$ws.Range("J29").Value = 2005000
$ws.Range("K29").Value = 120000
$ws.Range("L29").Value = 2005000
$ws.Range("M29").Value = -119710
$ws.Range("N29").Value = -2005580
$ws.Range("H50").Value = 20000
$ws.Range("J50").Value = 20000
$ws.Range("L50").Value = 20000
$ws.Range("N50").Value = -21262
$ws.Range("H53").Value = 0
$ws.Range("J53").Value = 0
$ws.Range("L53").Value = 0
$ws.Range("N53").ClearContents()
$ws.Range("H81").Value = 4267.8667
$ws.Range("I81").Value = 2555.5
$ws.Range("J81").Value = 4890.5454
$ws.Range("K81").Value = 5111
$ws.Range("L81").Value = 9781.0908
$ws.Range("M81").Value = -4050
$ws.Range("N81").Value = -11903.0908
$ws.Range("H84").Value = 4267.8667
$ws.Range("I84").Value = 2555.5
$ws.Range("J84").Value = 4890.5454
$ws.Range("K84").Value = 25555
$ws.Range("L84").Value = 48905.454
$ws.Range("M84").Value = -20251
$ws.Range("N84").Value = -59513.454
$ws.Range("H132").Value = 10001.857
$ws.Range("I132").Value = 9201.6
$ws.Range("J132").Value = 12002.5
$ws.Range("K132").Value = 27604.8
$ws.Range("L132").Value = 36007.5
$ws.Range("M132").Value = -25074.8
$ws.Range("N132").Value = -41067.5

